# Applies the content edits described by the commit:
#   "Some spelling and expansion of document."
#
# The underlying diff is dominated by Word's proofing engine splitting
# runs and inserting <w:proofErr/> markers around the same visible text
# (an artifact of a spell/grammar-check pass, not content available via
# the Word object model). The genuine content changes are:
#   1. "Our problem domain includes" -> "The problem domain includes"
#   2. "Testing should be silent, automatic and continuous" ->
#      "Testing should be automatic and continuous"   (drop "silent, ")
#   3. The "Incremental" paragraph gains a new trailing run of text:
#      "   9454948223"

$d = $word.ActiveDocument

# 1. "Our" -> "The" at the start of the problem-domain paragraph.
$d.Content.Find.Execute(
    "Our problem domain includes", $true, $false, $false, $false, $false,
    $true, 1, $false, "The problem domain includes", 2) | Out-Null

# 2. Drop "silent, " from the testing principle bullet.
$d.Content.Find.Execute(
    "Testing should be silent, automatic and continuous", $true, $false, $false, $false, $false,
    $true, 1, $false, "Testing should be automatic and continuous", 2) | Out-Null

# 3. Append a new run of digits after the word "Incremental".
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Incremental") {
        $insertPoint = $d.Range($p.Range.End - 1, $p.Range.End - 1)
        $insertPoint.InsertAfter("   9454948223") | Out-Null
        break
    }
}
